# Update the "取得日時" (retrieved datetime) timestamps on the ランサーズ sheet
# from 2025-11-26 18:22:19 to 2025-11-26 18:30:55 for rows 2 through 15 (column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-11-26 18:22:19"
$newValue = "2025-11-26 18:30:55"

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
